$d = $word.ActiveDocument

# Each of these paragraphs holds exactly one run, so replacing the whole
# paragraph range's text (rather than Find.Execute) keeps the existing
# (empty) run-properties element intact.

# 1. "Overhauled a web application ..." -> "Overhauled a Python web app ..."
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]10) -eq "Overhauled a web application used for filing and printing medical orders") {
        $p.Range.Text = "Overhauled a Python web app used for filing and printing medical orders"
        break
    }
}

# 2. Projects wiki link -> github.io wiki link
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]10) -eq "See: https://github.com/westurner/wiki/wiki/projects") {
        $p.Range.Text = "See: https://westurner.github.io/wiki/projects"
        break
    }
}

# 3. Contributions wiki link -> github.io wiki link
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]10) -eq "See: https://github.com/westurner/wiki/wiki/contributions") {
        $p.Range.Text = "See: https://westurner.github.io/wiki/contributions"
        break
    }
}

# 4. Add a new "https://westurner.github.io/dotfiles/" paragraph just before the
#    existing "https://github.com/westurner/dotfiles" paragraph (same style99 style).
#    NOTE: InsertParagraphBefore() reseats $p itself onto the freshly inserted
#    (empty) paragraph, so we can set its text directly afterwards.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]10) -eq "https://github.com/westurner/dotfiles") {
        $p.Range.InsertParagraphBefore() | Out-Null
        $p.Range.Text = "https://westurner.github.io/dotfiles/"
        break
    }
}

# 5. The "Python package with various Paver tasks" bullet becomes
#    "Python package with documentation" AND moves to appear right before the
#    "Configuration set for Bash, ZSH, Python, IPython, I3WM" bullet.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]10) -eq "Configuration set for Bash, ZSH, Python, IPython, I3WM") {
        $p.Range.InsertParagraphBefore() | Out-Null
        $p.Range.Text = "Python package with documentation"
        break
    }
}
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]10) -eq "Python package with various Paver tasks") {
        $p.Range.Delete() | Out-Null
        break
    }
}
